$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 156367.83
$ws.Range("I15").Value = 156367.83
$ws.Range("K15").Value = 469103.49
$ws.Range("M15").Value = -468934.49

$ws.Range("H33").Value = 112.61539
$ws.Range("I33").Value = 105.36364
$ws.Range("K33").Value = 105.36364
$ws.Range("M33").Value = 123.63636

$ws.Range("H62").Value = 4287470
$ws.Range("I62").Value = 6183151
$ws.Range("J62").Value = 22187.5
$ws.Range("K62").Value = 6183151
$ws.Range("L62").Value = 22187.5
$ws.Range("M62").Value = -6182527
$ws.Range("N62").Value = -23435.5

$ws.Range("H65").Value = 4287470
$ws.Range("I65").Value = 6183151
$ws.Range("J65").Value = 22187.5
$ws.Range("K65").Value = 30915755
$ws.Range("L65").Value = 110937.5
$ws.Range("M65").Value = -30912635
$ws.Range("N65").Value = -117177.5

$ws.Range("H107").Value = 358817.66
$ws.Range("I107").Value = 444812.16
$ws.Range("J107").Value = 507.16666
$ws.Range("K107").Value = 444812.16
$ws.Range("L107").Value = 507.16666
$ws.Range("M107").Value = -442892.16
$ws.Range("N107").Value = -4347.16666

$ws.Range("H132").Value = 423859.47
$ws.Range("I132").Value = 468708.44
$ws.Range("J132").Value = 35168.668
$ws.Range("K132").Value = 1406125.32
$ws.Range("L132").Value = 105506.004
$ws.Range("M132").Value = -1403595.32
$ws.Range("N132").Value = -110566.004

$ws.Range("H137").Value = 17544718
$ws.Range("I137").Value = 21739778
$ws.Range("J137").Value = 1734
$ws.Range("K137").Value = 65219334
$ws.Range("L137").Value = 5202
$ws.Range("M137").Value = -65216784
$ws.Range("N137").Value = -10302

$ws.Range("H138").Value = 1130.14
$ws.Range("I138").Value = 608.43396
$ws.Range("J138").Value = 1718.4468
$ws.Range("K138").Value = 1825.30188
$ws.Range("L138").Value = 5155.3404
$ws.Range("M138").Value = 3314.69812
$ws.Range("N138").Value = -15435.3404

$ws.Range("H141").Value = 2377.473
$ws.Range("I141").Value = 1472.9656
$ws.Range("J141").Value = 5656.3125
$ws.Range("K141").Value = 4418.8968
$ws.Range("L141").Value = 16968.9375
$ws.Range("M141").Value = 761.1031999999996
$ws.Range("N141").Value = -27328.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 49607.145
$ws.Range("I2").Value = 85001.164
$ws.Range("J2").Value = 2415.111
$ws.Range("K2").Value = 85001.164
$ws.Range("L2").Value = 2415.111
$ws.Range("M2").Value = -84888.164
$ws.Range("N2").Value = -2641.111

$ws.Range("H32").Value = 18880.047
$ws.Range("I32").Value = 3104.6833
$ws.Range("K32").Value = 3104.6833
$ws.Range("M32").Value = -2817.6833

$ws.Range("H61").Value = 1758.2449
$ws.Range("I61").Value = 1322.2273
$ws.Range("J61").Value = 5595.2
$ws.Range("K61").Value = 1322.2273
$ws.Range("L61").Value = 5595.2
$ws.Range("M61").Value = -1110.2273
$ws.Range("N61").Value = -6019.2

$ws.Range("H74").Value = 4810.75
$ws.Range("I74").Value = 1531.0667
$ws.Range("J74").Value = 14649.8
$ws.Range("K74").Value = 1531.0667
$ws.Range("L74").Value = 14649.8
$ws.Range("M74").Value = -657.0667000000001
$ws.Range("N74").Value = -16397.8

$ws.Range("H77").Value = 4810.75
$ws.Range("I77").Value = 1531.0667
$ws.Range("J77").Value = 14649.8
$ws.Range("K77").Value = 7655.333500000001
$ws.Range("L77").Value = 73249
$ws.Range("M77").Value = -3287.333500000001
$ws.Range("N77").Value = -81985

$ws.Range("H97").Value = 5868.8945
$ws.Range("I97").Value = 6843.625
$ws.Range("J97").Value = 670.3333
$ws.Range("K97").Value = 6843.625
$ws.Range("L97").Value = 670.3333
$ws.Range("M97").Value = -6347.625
$ws.Range("N97").Value = -1662.3333

$ws.Range("H116").Value = 49607.145
$ws.Range("I116").Value = 85001.164
$ws.Range("J116").Value = 2415.111
$ws.Range("K116").Value = 85001.164
$ws.Range("L116").Value = 2415.111
$ws.Range("M116").Value = -82707.164
$ws.Range("N116").Value = -7003.111

$ws.Range("H132").Value = 3625.75
$ws.Range("I132").Value = 3523.5789
$ws.Range("J132").Value = 3841.4443
$ws.Range("K132").Value = 10570.7367
$ws.Range("L132").Value = 11524.3329
$ws.Range("M132").Value = -8040.736699999999
$ws.Range("N132").Value = -16584.3329

$ws.Range("H136").Value = 1758.2449
$ws.Range("I136").Value = 1322.2273
$ws.Range("J136").Value = 5595.2
$ws.Range("K136").Value = 3966.6819
$ws.Range("L136").Value = 16785.6
$ws.Range("M136").Value = -1416.6819
$ws.Range("N136").Value = -21885.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 49607.145
$ws.Range("I3").Value = 85001.164
$ws.Range("J3").Value = 2415.111
$ws.Range("K3").Value = 85001.164
$ws.Range("L3").Value = 2415.111
$ws.Range("M3").Value = -84887.164
$ws.Range("N3").Value = -2643.111

$ws.Range("H94").Value = 1132.8
$ws.Range("I94").Value = 1045.8518
$ws.Range("J94").Value = 1426.25
$ws.Range("K94").Value = 1045.8518
$ws.Range("L94").Value = 1426.25
$ws.Range("M94").Value = -594.8517999999999
$ws.Range("N94").Value = -2328.25

$ws.Range("H105").Value = 268396.34
$ws.Range("I105").Value = 6001.037
$ws.Range("J105").Value = 912457.5600000001
$ws.Range("K105").Value = 6001.037
$ws.Range("L105").Value = 912457.5600000001
$ws.Range("M105").Value = -4254.037
$ws.Range("N105").Value = -915951.5600000001

$ws.Range("H134").Value = 2246.6724
$ws.Range("I134").Value = 1608.1459
$ws.Range("J134").Value = 5311.6
$ws.Range("K134").Value = 4824.4377
$ws.Range("L134").Value = 15934.8
$ws.Range("M134").Value = -2289.4377
$ws.Range("N134").Value = -21004.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1615.925
$ws.Range("I58").Value = 1013.0323
$ws.Range("J58").Value = 3692.5557
$ws.Range("K58").Value = 1013.0323
$ws.Range("L58").Value = 3692.5557
$ws.Range("M58").Value = -810.0323
$ws.Range("N58").Value = -4098.5557

$ws.Range("H62").Value = 33570.715
$ws.Range("I62").Value = 37499.168
$ws.Range("K62").Value = 37499.168
$ws.Range("M62").Value = -36875.168

$ws.Range("H65").Value = 33570.715
$ws.Range("I65").Value = 37499.168
$ws.Range("K65").Value = 187495.84
$ws.Range("M65").Value = -184375.84

$ws.Range("H119").Value = 42244.4
$ws.Range("J119").Value = 42244.4
$ws.Range("L119").Value = 42244.4
$ws.Range("N119").Value = -51920.4

$ws.Range("H136").Value = 1615.925
$ws.Range("I136").Value = 1013.0323
$ws.Range("J136").Value = 3692.5557
$ws.Range("K136").Value = 3039.0969
$ws.Range("L136").Value = 11077.6671
$ws.Range("M136").Value = -489.0969
$ws.Range("N136").Value = -16177.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 920
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = -2827

$ws.Range("H116").Value = 1719.75
$ws.Range("I116").Value = 959.6667
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 2879.0001
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = 562.9998999999998
$ws.Range("N116").Value = -18884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4000
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 7500
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 7500
$ws.Range("M29").Value = -210
$ws.Range("N29").Value = -8080

$ws.Range("H33").Value = 38981.816
$ws.Range("J33").Value = 38981.816
$ws.Range("L33").Value = 38981.816
$ws.Range("N33").Value = -39485.816

$ws.Range("H36").Value = 31500
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 31500
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 31500
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -32470

$ws.Range("H132").Value = 3270.7673
$ws.Range("I132").Value = 3096.7
$ws.Range("J132").Value = 3672.4614
$ws.Range("K132").Value = 9290.099999999999
$ws.Range("L132").Value = 11017.3842
$ws.Range("M132").Value = -6760.099999999999
$ws.Range("N132").Value = -16077.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1953.7693
$ws.Range("I68").Value = 1424.875
$ws.Range("K68").Value = 1424.875
$ws.Range("M68").Value = -675.875

$ws.Range("H70").Value = 22000
$ws.Range("J70").Value = 22000
$ws.Range("L70").Value = 22000
$ws.Range("N70").Value = -22540

$ws.Range("H71").Value = 1953.7693
$ws.Range("I71").Value = 1424.875
$ws.Range("K71").Value = 7124.375
$ws.Range("M71").Value = -3380.375

$ws.Range("H73").Value = 22000
$ws.Range("J73").Value = 22000
$ws.Range("L73").Value = 22000
$ws.Range("N73").Value = -23872

$ws.Range("H93").Value = 1251.5
$ws.Range("I93").Value = 503
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 503
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 745
$ws.Range("N93").Value = -4496

$ws.Range("H115").Value = 27334.285
$ws.Range("J115").Value = 27334.285
$ws.Range("L115").Value = 27334.285
$ws.Range("N115").Value = -29684.285

$ws.Range("H132").Value = 6892.8438
$ws.Range("I132").Value = 8980.134
$ws.Range("J132").Value = 5051.1177
$ws.Range("K132").Value = 26940.402
$ws.Range("L132").Value = 15153.3531
$ws.Range("M132").Value = -24410.402
$ws.Range("N132").Value = -20213.3531

$ws.Range("H136").Value = 4362.8
$ws.Range("I136").Value = 2375.1785
$ws.Range("J136").Value = 9000.583000000001
$ws.Range("K136").Value = 7125.5355
$ws.Range("L136").Value = 27001.749
$ws.Range("M136").Value = -4575.5355
$ws.Range("N136").Value = -32101.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 12800
$ws.Range("J19").Value = 3500
$ws.Range("L19").Value = 3500
$ws.Range("N19").Value = -3848

$ws.Range("H96").Value = 333334660
$ws.Range("I96").Value = 500001000
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 500001000
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -499999627
$ws.Range("N96").Value = -4746

$ws.Range("H122").Value = 84356.25
$ws.Range("I122").Value = 111932.78
$ws.Range("J122").Value = 1626.6666
$ws.Range("K122").Value = 335798.34
$ws.Range("L122").Value = 4879.9998
$ws.Range("M122").Value = -333348.34
$ws.Range("N122").Value = -9779.9998

$ws.Range("H126").Value = 56734.89
$ws.Range("I126").Value = 67735.2
$ws.Range("J126").Value = 1733.3334
$ws.Range("K126").Value = 203205.6
$ws.Range("L126").Value = 5200.0002
$ws.Range("M126").Value = -200735.6
$ws.Range("N126").Value = -10140.0002

$ws.Range("H132").Value = 8930731
$ws.Range("I132").Value = 13159947
$ws.Range("J132").Value = 2384.4443
$ws.Range("K132").Value = 39479841
$ws.Range("L132").Value = 7153.3329
$ws.Range("M132").Value = -39477311
$ws.Range("N132").Value = -12213.3329

$ws.Range("H136").Value = 19958.754
$ws.Range("I136").Value = 23254.182
$ws.Range("J136").Value = 3847.7778
$ws.Range("K136").Value = 69762.546
$ws.Range("L136").Value = 11543.3334
$ws.Range("M136").Value = -67212.546
$ws.Range("N136").Value = -16643.3334
